$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.151.44"
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").Value = "1.849.83"
$ws.Range("E3").Value = "  -0.80%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Value = "'235.41"
$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("D7").Value = "'0.4724"
$ws.Range("E7").Value = "  +1.24%  "

$ws.Range("D8").Value = "'0.2905"
$ws.Range("E8").Value = "  +2.51%  "

$ws.Range("D9").Value = "'0.06533"
$ws.Range("E9").Value = "  +0.25%  "

$ws.Range("D10").Value = "'21.58"
$ws.Range("E10").Value = "  +0.43%  "

$ws.Range("D11").Value = "'0.07938"
$ws.Range("E11").Value = "  +1.10%  "

$ws.Range("D12").Value = "'97.54"
$ws.Range("E12").Value = "  +0.17%  "

$ws.Range("D13").Value = "1.858.71"
$ws.Range("E13").Value = "  -0.71%  "

$ws.Range("D14").Value = "'5.081"
$ws.Range("E14").Value = "  -0.26%  "

$ws.Range("D15").Value = "'0.6742"
$ws.Range("E15").Value = "  +0.29%  "

$ws.Range("D16").Value = "'267.81"
$ws.Range("E16").Value = "  -4.51%  "

$ws.Range("D17").Value = "30.138.31"
$ws.Range("E17").Value = "  -0.32%  "

$ws.Range("D18").Value = "'13.57"
$ws.Range("E18").Value = "  +6.99%  "

$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = "  +0.20%  "

$ws.Range("D20").Value = "'0.000007543"
$ws.Range("E20").Value = "  +3.65%  "

$ws.Range("D21").Value = "2.103.72"
$ws.Range("E21").Value = "  -0.39%  "

$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("D23").Value = "'5.227"
$ws.Range("E23").Value = "  -4.65%  "

$ws.Range("D24").Value = "'6.125"
$ws.Range("E24").Value = "  -0.32%  "

$ws.Range("D25").Value = "'165.84"
$ws.Range("E25").Value = "  +1.03%  "

$ws.Range("D26").Value = "'9.134"
$ws.Range("E26").Value = "  -0.39%  "

$ws.Range("D27").Value = "'18.76"
$ws.Range("E27").Value = "  -1.77%  "

$ws.Range("D28").Value = "'1.921"
$ws.Range("E28").Value = "  -0.37%  "

$ws.Range("D29").Value = "'1.395"
$ws.Range("E29").Value = "  +1.07%  "

$ws.Range("D30").Value = "'0.09860"
$ws.Range("E30").Value = "  +1.79%  "

$ws.Range("D31").Value = "'1.463"
$ws.Range("E31").Value = "  -0.94%  "

$ws.Range("D32").Value = "'4.269"
$ws.Range("E32").Value = "  -3.02%  "

$ws.Range("D33").Value = "'3.994"
$ws.Range("E33").Value = "  -2.39%  "

$ws.Range("D34").Value = "'0.04673"
$ws.Range("E34").Value = "  -0.26%  "

$ws.Range("D35").Value = "'1.117"
$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("D36").Value = "'0.6960"
$ws.Range("E36").Value = "  -1.34%  "

$ws.Range("D37").Value = "'2.710"
$ws.Range("E37").Value = "  -0.69%  "

$ws.Range("D38").Value = "'0.01862"
$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("D39").Value = "'2.606"
$ws.Range("E39").Value = "  +2.88%  "

$ws.Range("D40").Value = "'6.310"
$ws.Range("E40").Value = "  +1.06%  "

$ws.Range("D41").Value = "'73.13"
$ws.Range("E41").Value = "  +0.08%  "

$ws.Range("D42").Value = "'1.927"
$ws.Range("E42").Value = "  -0.70%  "

$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("D44").Value = "'0.8351"
$ws.Range("E44").Value = "  -1.44%  "

$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'102.97"
$ws.Range("E45").Value = "  -0.89%  "

$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").Value = "'0.4117"
$ws.Range("E46").Value = "  -1.24%  "

$ws.Range("D47").Value = "'939.82"
$ws.Range("E47").Value = "  +0.41%  "

$ws.Range("D48").Value = "'9.077"
$ws.Range("E48").Value = "  -0.93%  "

$ws.Range("D49").Value = "'6.947"
$ws.Range("E49").Value = "  -3.38%  "

$ws.Range("D50").Value = "'33.61"
$ws.Range("E50").Value = "  -1.50%  "

$ws.Range("D51").Value = "'0.05650"
$ws.Range("E51").Value = "  +0.34%  "
